# Sync attendance_reports: normalize "Recorded By" (column G) ordering.
# Several rows had their comma-separated list of recorders written in the
# opposite order to the canonical form used elsewhere in the workbook.
# Re-order those exact values back to the canonical order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of "out of order" text -> canonical re-ordered text.
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, system, System";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "admin@admin.com, System"             = "System, admin@admin.com";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$lastRow = $ws.UsedRange.Rows.Count
$col = 7   # Column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
